$d = $word.ActiveDocument

function Find-ParagraphByText($needle) {
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text -match [regex]::Escape($needle)) {
            return $p
        }
    }
    return $null
}

function Append-RedText($paragraph, $text) {
    $endPos = $paragraph.Range.End - 1
    $insertPoint = $d.Range($endPos, $endPos)
    $insertPoint.InsertAfter($text)
    $colored = $d.Range($endPos, $endPos + $text.Length)
    $colored.Font.Color = 255
}

# --- Step 1: remove the "-Liste des musiques en attentes dans la playlist" paragraph entirely ---
$pList = Find-ParagraphByText("des musiques en attentes dans la playlist")
$pList.Range.Delete()

# --- Step 2: append to "-Tableau associatif..." paragraph ---
$pTab = Find-ParagraphByText("Tableau associatif")
Append-RedText $pTab ", trié dans l’ordre de la playlist"
Append-RedText $pTab ". La première musique est celle en cours"

# --- Step 3: remove "-Nom de la musique en cours" paragraph entirely ---
$pNom = Find-ParagraphByText("Nom de la musique en cours")
$pNom.Range.Delete()

# --- Step 4: append " (s)" to "-Timestamp de la musique en cours" ---
$pTime = Find-ParagraphByText("Timestamp de la musique en cours")
Append-RedText $pTime " (s)"

# --- Step 5: append " (" and "s)" to "-Durée de la musique en cours" and add bookmark ---
$pDur = Find-ParagraphByText("Durée de la musique en cours")
Append-RedText $pDur " ("
Append-RedText $pDur "s)"
Write-Output "pDur range:"
Write-Output $pDur.Range.Start
Write-Output $pDur.Range.End
Write-Output $pDur.Range.Text
$d.Bookmarks.Add("TestBM", $pDur.Range)
Write-Output "bookmark added, exists:"
Write-Output $d.Bookmarks.Exists("TestBM")
